# Applies the odds updates described in the commit diff to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("T2").Value = 1.55
$ws.Range("AD2").Value = 19.5
$ws.Range("AI2").Value = 50
$ws.Range("AL2").Value = 29
$ws.Range("AO2").Value = 36

# Row 3
$ws.Range("F3").Value = 1.92
$ws.Range("H3").Value = 4.1
$ws.Range("O3").Value = 1.23
$ws.Range("R3").Value = 1.53
$ws.Range("AC3").Value = 9.800000000000001

# Row 4
$ws.Range("L4").Value = 1.01
$ws.Range("N4").Value = 4.2

# Row 6
$ws.Range("F6").Value = 11.5
$ws.Range("G6").Value = 19.5
$ws.Range("H6").Value = 1.2

# Row 8
$ws.Range("H8").Value = 3.55
$ws.Range("I8").Value = 3.6
$ws.Range("O8").Value = 1.37
$ws.Range("U8").Value = 2.08
$ws.Range("AB8").Value = 9.199999999999999
$ws.Range("AF8").Value = 13.5
$ws.Range("AH8").Value = 19
$ws.Range("AJ8").Value = 30
$ws.Range("AM8").Value = 110

# Row 9
$ws.Range("F9").Value = 1.58
$ws.Range("G9").Value = 1.59
$ws.Range("N9").Value = 5.2
$ws.Range("O9").Value = 1.22
$ws.Range("P9").Value = 2.5
$ws.Range("Q9").Value = 1.65
$ws.Range("R9").Value = 1.56
$ws.Range("S9").Value = 2.68
$ws.Range("T9").Value = 1.77
$ws.Range("U9").Value = 2.26
$ws.Range("V9").Value = 1.18
$ws.Range("W9").Value = 2.7
$ws.Range("X9").Value = 24
$ws.Range("Y9").Value = 27
$ws.Range("AG9").Value = 9.4
$ws.Range("AH9").Value = 19
$ws.Range("AJ9").Value = 15
$ws.Range("AM9").Value = 85

# Row 10
$ws.Range("I10").Value = 4.5
$ws.Range("J10").Value = 3.5
$ws.Range("V10").Value = 1.28

# Row 11
$ws.Range("F11").Value = 2.32
$ws.Range("G11").Value = 2.34
$ws.Range("P11").Value = 1.93
$ws.Range("T11").Value = 1.82
$ws.Range("U11").Value = 2.14
$ws.Range("W11").Value = 1.74
$ws.Range("AB11").Value = 9.800000000000001
$ws.Range("AN11").Value = 19.5

# Row 12
$ws.Range("F12").Value = 8.4
$ws.Range("G12").Value = 8.6
$ws.Range("H12").Value = 1.44
$ws.Range("I12").Value = 1.45
$ws.Range("J12").Value = 5.2
$ws.Range("K12").Value = 5.3
$ws.Range("T12").Value = 1.93
$ws.Range("V12").Value = 3.2
$ws.Range("X12").Value = 22
$ws.Range("Y12").Value = 9.4
$ws.Range("AA12").Value = 12
$ws.Range("AE12").Value = 14.5
$ws.Range("AF12").Value = 75
$ws.Range("AG12").Value = 30
$ws.Range("AH12").Value = 24
$ws.Range("AJ12").Value = 250
$ws.Range("AK12").Value = 120
$ws.Range("AM12").Value = 130
$ws.Range("AO12").Value = 6.2
